$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version
$ws.Range("B3").Value = "0.1.7"

# Update Status
$ws.Range("B6").Value = "draft"

# Update Date
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"

# Update the two Contact rows (10 and 11) with the new publisher/contact info
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new "Jurisdiction" row after the Contact rows (before Description),
# copying formatting from the row above so the new row matches existing style.
$ws.Rows.Item(12).Insert()
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
